$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking values
# are not auto-converted to numbers by Excel (matches original inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "24.414.01"
$ws.Cells.Item(2, 5).Value = "  +8.60%  "
$ws.Cells.Item(3, 4).Value = "1.680.19"
$ws.Cells.Item(3, 5).Value = "  +4.15%  "
$ws.Cells.Item(4, 4).Value = "1.006"
$ws.Cells.Item(4, 5).Value = "  +0.34%  "
$ws.Cells.Item(5, 4).Value = "307.66"
$ws.Cells.Item(5, 5).Value = "  +0.90%  "
$ws.Cells.Item(6, 4).Value = "0.9999"
$ws.Cells.Item(6, 5).Value = "  +0.93%  "
$ws.Cells.Item(7, 4).Value = "0.3713"
$ws.Cells.Item(7, 5).Value = "  +0.60%  "
$ws.Cells.Item(8, 4).Value = "0.3441"
$ws.Cells.Item(8, 5).Value = "  +0.62%  "
$ws.Cells.Item(9, 4).Value = "47.74"
$ws.Cells.Item(9, 5).Value = "  +12.72%  "
$ws.Cells.Item(10, 4).Value = "1.172"
$ws.Cells.Item(10, 5).Value = "  +2.81%  "
$ws.Cells.Item(11, 4).Value = "0.07269"
$ws.Cells.Item(11, 5).Value = "  +2.54%  "
$ws.Cells.Item(12, 4).Value = "1.003"
$ws.Cells.Item(12, 5).Value = "  +0.33%  "
$ws.Cells.Item(13, 4).Value = "6.114"
$ws.Cells.Item(13, 5).Value = "  +2.95%  "
$ws.Cells.Item(14, 4).Value = "20.26"
$ws.Cells.Item(14, 5).Value = "  +2.18%  "
$ws.Cells.Item(15, 4).Value = "6.744"
$ws.Cells.Item(15, 5).Value = "  +1.10%  "
$ws.Cells.Item(16, 4).Value = "1.680.23"
$ws.Cells.Item(16, 5).Value = "  +4.63%  "
$ws.Cells.Item(17, 5).Value = "  +1.62%  "
$ws.Cells.Item(18, 4).Value = "1.0000"
$ws.Cells.Item(18, 5).Value = "  +0.90%  "
$ws.Cells.Item(19, 4).Value = "0.06683"
$ws.Cells.Item(19, 5).Value = "  -1.34%  "
$ws.Cells.Item(20, 4).Value = "81.13"
$ws.Cells.Item(20, 5).Value = "  +3.47%  "
$ws.Cells.Item(21, 4).Value = "16.46"
$ws.Cells.Item(21, 5).Value = "  +2.01%  "
$ws.Cells.Item(22, 4).Value = "6.123"
$ws.Cells.Item(22, 5).Value = "  +1.19%  "
$ws.Cells.Item(23, 5).Value = "  +2.50%  "
$ws.Cells.Item(24, 4).Value = "24.358.26"
$ws.Cells.Item(24, 5).Value = "  +8.37%  "
$ws.Cells.Item(25, 4).Value = "2.458"
$ws.Cells.Item(25, 5).Value = "  +2.94%  "
$ws.Cells.Item(26, 2).Value = "LEO"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(26, 4).Value = "3.356"
$ws.Cells.Item(26, 5).Value = "  -13.36%  "
$ws.Cells.Item(27, 2).Value = "LidoDAOToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(27, 4).Value = "2.661"
$ws.Cells.Item(27, 5).Value = "  +3.91%  "
$ws.Cells.Item(28, 2).Value = "Monero"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(28, 4).Value = "153.56"
$ws.Cells.Item(28, 5).Value = "  +2.24%  "
$ws.Cells.Item(29, 2).Value = "EthereumClassic"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(29, 4).Value = "19.53"
$ws.Cells.Item(29, 5).Value = "  -0.45%  "
$ws.Cells.Item(30, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(30, 4).Value = "1.862.86"
$ws.Cells.Item(30, 5).Value = "  +4.19%  "
$ws.Cells.Item(31, 2).Value = "BitcoinCash"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(31, 4).Value = "127.62"
$ws.Cells.Item(31, 5).Value = "  +3.70%  "
$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).Value = "6.314"
$ws.Cells.Item(32, 5).Value = "  +1.86%  "
$ws.Cells.Item(33, 2).Value = "HuobiToken"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(33, 4).Value = "4.062"
$ws.Cells.Item(33, 5).Value = "  +0.36%  "
$ws.Cells.Item(34, 2).Value = "ImmutableX"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(34, 4).Value = "0.9745"
$ws.Cells.Item(34, 5).Value = "  +1.91%  "
$ws.Cells.Item(35, 2).Value = "Stellar"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(35, 4).Value = "0.08448"
$ws.Cells.Item(35, 5).Value = "  +2.31%  "
$ws.Cells.Item(36, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(36, 4).Value = "1.698"
$ws.Cells.Item(36, 5).Value = "  +2.11%  "
$ws.Cells.Item(37, 2).Value = "Aptos"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(37, 4).Value = "12.36"
$ws.Cells.Item(37, 5).Value = "  +2.71%  "
$ws.Cells.Item(38, 2).Value = "Hedera"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(38, 4).Value = "0.06485"
$ws.Cells.Item(38, 5).Value = "  +5.86%  "
$ws.Cells.Item(39, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(39, 4).Value = "5.353"
$ws.Cells.Item(39, 5).Value = "  +1.40%  "
$ws.Cells.Item(40, 2).Value = "FraxShare"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(40, 4).Value = "8.895"
$ws.Cells.Item(40, 5).Value = "  +3.15%  "
$ws.Cells.Item(41, 2).Value = "VeChain"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(41, 4).Value = "0.02329"
$ws.Cells.Item(41, 5).Value = "  +4.17%  "
$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).Value = "1.254"
$ws.Cells.Item(42, 5).Value = "  -1.59%  "
$ws.Cells.Item(43, 2).Value = "Algorand"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(43, 4).Value = "0.2099"
$ws.Cells.Item(43, 5).Value = "  +3.49%  "
$ws.Cells.Item(44, 2).Value = "TheSandbox"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(44, 4).Value = "0.6162"
$ws.Cells.Item(44, 5).Value = "  +3.70%  "
$ws.Cells.Item(45, 2).Value = "Frax"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(45, 4).Value = "0.9998"
$ws.Cells.Item(45, 5).Value = "  +0.85%  "
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).Value = "13.18"
$ws.Cells.Item(46, 5).Value = "  -0.13%  "
$ws.Cells.Item(47, 2).Value = "PancakeSwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(47, 4).Value = "3.777"
$ws.Cells.Item(47, 5).Value = "  -1.42%  "
$ws.Cells.Item(48, 2).Value = "Decentraland"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(48, 4).Value = "0.5938"
$ws.Cells.Item(48, 5).Value = "  +3.78%  "
$ws.Cells.Item(49, 2).Value = "Quant"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(49, 4).Value = "127.39"
$ws.Cells.Item(49, 5).Value = "  +0.06%  "
$ws.Cells.Item(50, 2).Value = "NEARProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(50, 4).Value = "2.027"
$ws.Cells.Item(50, 5).Value = "  +2.12%  "
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "0.07184"
$ws.Cells.Item(51, 5).Value = "  +5.25%  "
